# Update crypto price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.101.82"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.990.98"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.12"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.52"
$ws.Range("E6").Value = "  -4.59%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.627"
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.23"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  -4.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.37"
$ws.Range("E13").Value = "  -3.55%  "
$ws.Range("D14").Value = "3.461.43"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.67"
$ws.Range("E15").Value = "  -3.18%  "
$ws.Range("D16").Value = "2.988.34"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "52.135.50"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.47"
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  -5.61%  "
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.47"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "263.91"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.89"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("E31").Value = "  -3.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.62"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.90"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("E39").Value = "  -5.57%  "
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.69"
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.89"
$ws.Range("E44").Value = "  +8.12%  "
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").Value = "2.125.16"
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("E48").Value = "  -5.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.242"
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("E51").Value = "  -1.21%  "
